# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.949.52"
$ws.Range("E2").Value = "  -1.87%  "

$ws.Range("D3").Value = "1.635.61"
$ws.Range("E3").Value = "  -1.96%  "

$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = "  +0.56%  "

$ws.Range("D5").Value = "'215.18"
$ws.Range("E5").Value = "  -1.66%  "

$ws.Range("D6").Value = "'0.5020"
$ws.Range("E6").Value = "  -2.41%  "

$ws.Range("D7").Value = "'1.013"
$ws.Range("E7").Value = "  +0.55%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.06410"
$ws.Range("E9").Value = "  -0.96%  "

$ws.Range("D10").Value = "'19.49"
$ws.Range("E10").Value = "  -2.64%  "

$ws.Range("E11").Value = "  +1.32%  "

$ws.Range("D12").Value = "1.639.68"
$ws.Range("E12").Value = "  -2.00%  "

$ws.Range("D13").Value = "'4.263"
$ws.Range("E13").Value = "  -1.94%  "

$ws.Range("D14").Value = "1.860.09"
$ws.Range("E14").Value = "  -2.03%  "

$ws.Range("D15").Value = "'0.5444"
$ws.Range("E15").Value = "  -2.23%  "

$ws.Range("D16").Value = "0.0₅7956"
$ws.Range("E16").Value = "  -1.21%  "

$ws.Range("D17").Value = "'63.48"
$ws.Range("E17").Value = "  -1.98%  "

$ws.Range("D18").Value = "25.935.69"
$ws.Range("E18").Value = "  -2.05%  "

$ws.Range("D19").Value = "'1.015"
$ws.Range("E19").Value = "  +0.90%  "

$ws.Range("D20").Value = "'204.80"
$ws.Range("E20").Value = "  -2.71%  "

$ws.Range("D21").Value = "'4.312"
$ws.Range("E21").Value = "  -2.78%  "

$ws.Range("D22").Value = "'9.997"
$ws.Range("E22").Value = "  -1.31%  "

$ws.Range("D23").Value = "'5.970"
$ws.Range("E23").Value = "  +1.28%  "

$ws.Range("D24").Value = "'1.014"
$ws.Range("E24").Value = "  +0.66%  "

$ws.Range("D25").Value = "'1.969"
$ws.Range("E25").Value = "  +13.42%  "

$ws.Range("D26").Value = "'141.85"
$ws.Range("E26").Value = "  -2.56%  "

$ws.Range("E27").Value = "  -1.14%  "

$ws.Range("D28").Value = "'15.76"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").Value = "'6.805"
$ws.Range("E29").Value = "  -2.98%  "

$ws.Range("D31").Value = "'0.04997"
$ws.Range("E31").Value = "  -4.33%  "

$ws.Range("D32").Value = "'3.263"
$ws.Range("E32").Value = "  -3.26%  "

$ws.Range("D33").Value = "'3.195"
$ws.Range("E33").Value = "  -0.82%  "

$ws.Range("D34").Value = "'1.539"
$ws.Range("E34").Value = "  -2.92%  "

$ws.Range("D35").Value = "'2.343"
$ws.Range("E35").Value = "  -1.40%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.8896"
$ws.Range("E36").Value = "  -3.75%  "

$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'2.614"
$ws.Range("E37").Value = "  -5.47%  "

$ws.Range("D38").Value = "'0.5652"
$ws.Range("E38").Value = "  -1.68%  "

$ws.Range("D39").Value = "1.114.41"
$ws.Range("E39").Value = "  -4.56%  "

$ws.Range("E40").Value = "  -2.61%  "

$ws.Range("D41").Value = "'2.581"
$ws.Range("E41").Value = "  +0.11%  "

$ws.Range("D42").Value = "'1.014"
$ws.Range("E42").Value = "  +0.71%  "

$ws.Range("D43").Value = "'5.607"
$ws.Range("E43").Value = "  -0.70%  "

$ws.Range("D44").Value = "'0.8167"
$ws.Range("E44").Value = "  -3.43%  "

$ws.Range("D45").Value = "'99.72"
$ws.Range("E45").Value = "  -0.60%  "

$ws.Range("D46").Value = "1.770.97"
$ws.Range("E46").Value = "  -2.11%  "

$ws.Range("E47").Value = "  +1.21%  "

$ws.Range("D48").Value = "'0.4546"
$ws.Range("E48").Value = "  +1.19%  "

$ws.Range("D49").Value = "'1.016"
$ws.Range("E49").Value = "  +1.02%  "

$ws.Range("D50").Value = "'54.83"
$ws.Range("E50").Value = "  -2.36%  "

$ws.Range("D51").Value = "'0.05035"
$ws.Range("E51").Value = "  -1.67%  "

